# Css changes in main page
# Update bug_reports sheet: change a handful of status/date cells and
# append four new bug rows (16-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "fixed"
$ws.Range("D2").Value = "Closed"
$ws.Range("F2").Value = "2023-08-11 14:40:28"

# --- Row 3 ---
$ws.Range("F3").Value = "2023-08-11 15:15:13"

# --- Row 7 ---
$ws.Range("D7").Value = "Closed"
$ws.Range("F7").Value = "2023-08-04 16:07:02"

# --- Row 9 ---
$ws.Range("D9").Value = "Closed"
$ws.Range("F9").Value = "2023-08-04 16:06:53"

# --- Row 15 ---
$ws.Range("D15").Value = "In Progress"
$ws.Range("F15").Value = "2023-08-04 08:24:28"

# --- New Row 16 ---
$ws.Range("A16").Value = 567
$ws.Range("B16").Value = "fdsfsdfs"
$ws.Range("C16").Value = "fdsfsdfs"
$ws.Range("D16").Value = "Closed"
$ws.Range("E16").Value = "2023-08-04 08:24:46"
$ws.Range("F16").Value = "2023-08-04 16:07:19"

# --- New Row 17 ---
$ws.Range("A17").Value = 3242
$ws.Range("B17").Value = "sewrewr"
$ws.Range("C17").Value = "sewrewr"
$ws.Range("D17").Value = "Closed"
$ws.Range("E17").Value = "2023-08-04 08:25:01"
$ws.Range("F17").Value = "2023-08-04 08:29:45"

# --- New Row 18 ---
$ws.Range("A18").Value = 314
$ws.Range("B18").Value = "erewrw"
$ws.Range("C18").Value = "erewrw"
$ws.Range("D18").Value = "Closed"
$ws.Range("E18").Value = "2023-08-04 08:25:11"
$ws.Range("F18").Value = "2023-08-04 16:07:29"

# --- New Row 19 ---
$ws.Range("A19").Value = 999
$ws.Range("B19").Value = "good work as bug"
$ws.Range("C19").Value = "good work as bug is working"
$ws.Range("D19").Value = "In Progress"
$ws.Range("E19").Value = "2023-08-11 15:47:32"
$ws.Range("F19").Value = "2023-08-11 15:48:45"

# Column A on the existing bug rows carries a bold/bordered/centered style
# (s="1"). Copy that formatting from the last original row (A15) onto the
# newly added rows so they match the same look, without creating
# duplicate style entries.
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)
